$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells (row 1)
$ws.Range('A1').Value = 'mx_state'
$ws.Range('B1').Value = 'mx_municipality'
$ws.Range('C1').Value = 'n_matriculas'
$ws.Range('D1').Value = 'pct_matriculas'

# Title-case the connector words (de, del, el, la, los, y) throughout
# the Estado/Municipio text columns
$updates = @{
    'B7' = 'Rincón De Romos'
    'B25' = 'Comitán De Domínguez'
    'B34' = 'Mazapa De Madero'
    'B55' = 'Hidalgo Del Parral'
    'A86' = 'Ciudad De México'
    'B90' = 'Cuajimalpa De Morelos'
    'B101' = 'Coneto De Comonfort'
    'A122' = 'Estado De México'
    'B122' = 'Acambay De Ruíz Castañeda'
    'B125' = 'Almoloya De Juárez'
    'B129' = 'Atizapán De Zaragoza'
    'B132' = 'Chapa De Mota'
    'B134' = 'Coacalco De Berriozábal'
    'B139' = 'Ecatepec De Morelos'
    'B140' = 'Ixtapan De La Sal'
    'B144' = 'Naucalpan De Juárez'
    'B149' = 'San Felipe Del Progreso'
    'B150' = 'Soyaniquilpan De Juárez'
    'B157' = 'Tlalnepantla De Baz'
    'B162' = 'Valle De Bravo'
    'B163' = 'Villa De Allende'
    'B170' = 'San Miguel De Allende'
    'B171' = 'Apaseo El Alto'
    'B175' = 'Dolores Hidalgo Cuna De La Independencia Nacional'
    'B179' = 'Jaral Del Progreso'
    'B185' = 'Purísima Del Rincón'
    'B190' = 'San Francisco Del Rincón'
    'B192' = 'San Luis De La Paz'
    'B193' = 'Santa Cruz De Juventino Rosas'
    'B195' = 'Silao De La Victoria'
    'B198' = 'Valle De Santiago'
    'B201' = 'Acapulco De Juárez'
    'B202' = 'Ajuchitlán Del Progreso'
    'B206' = 'Atoyac De Álvarez'
    'B208' = 'Chilapa De Álvarez'
    'B209' = 'Chilpancingo De Los Bravo'
    'B210' = 'Coahuayutla De José María Izazaga'
    'B217' = 'Cutzamala De Pinzón'
    'B220' = 'Huitzuco De Los Figueroa'
    'B221' = 'Iguala De La Independencia'
    'B222' = 'Ixcateopan De Cuauhtémoc'
    'B223' = 'Zihuatanejo De Azueta'
    'B232' = 'Taxco De Alarcón'
    'B233' = 'Técpan De Galeana'
    'B235' = 'Tepecoacuilco De Trujano'
    'B243' = 'Atotonilco El Grande'
    'B246' = 'Cuautepec De Hinojosa'
    'B250' = 'Mineral Del Chico'
    'B251' = 'Mineral Del Monte'
    'B253' = 'Nopala De Villagrán'
    'B254' = 'Pachuca De Soto'
    'B256' = 'Progreso De Obregón'
    'B257' = 'Santiago Tulantepec De Lugo Guerrero'
    'B259' = 'Tepeji Del Río De Ocampo'
    'B260' = 'Tezontepec De Aldama'
    'B262' = 'Tula De Allende'
    'B263' = 'Tulancingo De Bravo'
    'B264' = 'Villa De Tezontepec'
    'B265' = 'Zacualtipán De Ángeles'
    'B269' = 'Acatlán De Juárez'
    'B273' = 'Atotonilco El Alto'
    'B274' = 'Autlán De Navarro'
    'B279' = 'Cañadas De Obregón'
    'B292' = 'Huejuquilla El Alto'
    'B295' = 'Jilotlán De Los Dolores'
    'B299' = 'Lagos De Moreno'
    'B304' = 'Ojuelos De Jalisco'
    'B308' = 'San Juanito De Escobedo'
    'B309' = 'San Martín De Bolaños'
    'B310' = 'San Miguel El Alto'
    'B311' = 'San Sebastián Del Oeste'
    'B314' = 'Talpa De Allende'
    'B315' = 'Tamazula De Gordiano'
    'B319' = 'Teocuitatlán De Corona'
    'B320' = 'Tepatitlán De Morelos'
    'B322' = 'Tizapán El Alto'
    'B323' = 'Tlajomulco De Zúñiga'
    'B331' = 'Unión De Tula'
    'B332' = 'Valle De Guadalupe'
    'B335' = 'Yahualica De González Gallo'
    'B336' = 'Zacoalco De Torres'
    'B339' = 'Zapotitlán De Vadillo'
    'B340' = 'Zapotlán El Grande'
    'B361' = 'Coalcomán De Vázquez Pallares'
    'B433' = 'Coatlán Del Río'
    'B442' = 'Tetela Del Volcán'
    'B443' = 'Tlaltizapán De Zapata'
    'B452' = 'Amatlán De Cañas'
    'B453' = 'Bahía De Banderas'
    'B456' = 'Ixtlán Del Río'
    'B462' = 'Santa María Del Oro'
    'B475' = 'San Nicolás De Los Garza'
    'B477' = 'Acatlán De Pérez Figueroa'
    'B479' = 'Heroica Ciudad De Ejutla De Crespo'
    'B480' = 'Heroica Ciudad De Tlaxiaco'
    'B482' = 'Miahuatlán De Porfirio Díaz'
    'B483' = 'Oaxaca De Juárez'
    'B484' = 'Ocotlán De Morelos'
    'B485' = 'Putla Villa De Guerrero'
    'B503' = 'San Pedro El Alto'
    'B506' = 'San Pedro Y San Pablo Tequixtepec'
    'B508' = 'Santa Inés Del Monte'
    'B523' = 'Teococuilco De Marcos Pérez'
    'B524' = 'Teotitlán De Flores Magón'
    'B525' = 'Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca'
    'B526' = 'Tlacolula De Matamoros'
    'B527' = 'Villa Sola De Vega'
    'B528' = 'Zapotitlán Del Río'
    'B529' = 'Zimatlán De Álvarez'
    'B536' = 'Chalchicomula De Sesma'
    'B541' = 'Cuayuca De Andrade'
    'B545' = 'Huehuetlán El Grande'
    'B557' = 'Tepexi De Rodríguez'
    'B560' = 'Totoltepec De Guerrero'
    'B567' = 'Cadereyta De Montes'
    'B571' = 'Jalpan De Serra'
    'B573' = 'Pinal De Amoles'
    'B576' = 'San Juan Del Río'
    'B584' = 'Armadillo De Los Infante'
    'B585' = 'Axtla De Terrazas'
    'B588' = 'Ciudad Del Maíz'
    'B592' = 'Mexquitic De Carmona'
    'B597' = 'Villa De Ramos'
    'B598' = 'Villa De Reyes'
    'B640' = 'Jalpa De Méndez'
    'B664' = 'Alto Lucero De Gutiérrez Barrios'
    'B670' = 'Cosamaloapan De Carpio'
    'B675' = 'Ixhuatlán De Madero'
    'B676' = 'Ixhuatlán Del Sureste'
    'B679' = 'Juchique De Ferrer'
    'B682' = 'Lerdo De Tejada'
    'B687' = 'Ozuluama De Mascareñas'
    'B701' = 'Vega De Alatorre'
    'B723' = 'Jiménez Del Teul'
    'B728' = 'Moyahua De Estrada'
    'B729' = 'Nochistlán De Mejía'
    'B738' = 'Teúl De González Ortega'
    'B739' = 'Tlaltenango De Sánchez Román'
    'B741' = 'Villa De Cos'
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# Remove the trailing metadata/footer rows (749-753); row 748 was already blank
$ws.Range('A749:A753').EntireRow.Delete()

Write-Host "done"